$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 484.85
$ws.Range("J17").Value = 484.59323
$ws.Range("L17").Value = 1453.77969
$ws.Range("N17").Value = -1789.77969
$ws.Range("H57").Value = 20148
$ws.Range("J57").Value = 20148
$ws.Range("L57").Value = 60444
$ws.Range("N57").Value = -61442
$ws.Range("H62").Value = 218755380
$ws.Range("I62").Value = 100008600
$ws.Range("J62").Value = 416666660
$ws.Range("K62").Value = 100008600
$ws.Range("L62").Value = 416666660
$ws.Range("M62").Value = -100007976
$ws.Range("N62").Value = -416667908
$ws.Range("H65").Value = 218755380
$ws.Range("I65").Value = 100008600
$ws.Range("J65").Value = 416666660
$ws.Range("K65").Value = 500043000
$ws.Range("L65").Value = 2083333300
$ws.Range("M65").Value = -500039880
$ws.Range("N65").Value = -2083339540
$ws.Range("H98").Value = 26353014
$ws.Range("I98").Value = 14287568
$ws.Range("J98").Value = 39346572
$ws.Range("K98").Value = 14287568
$ws.Range("L98").Value = 39346572
$ws.Range("M98").Value = -14286070
$ws.Range("N98").Value = -39349568
$ws.Range("H112").Value = 1159029.5
$ws.Range("I112").Value = 900
$ws.Range("J112").Value = 1227154.9
$ws.Range("K112").Value = 2700
$ws.Range("L112").Value = 3681464.7
$ws.Range("M112").Value = -1592
$ws.Range("N112").Value = -3683680.7
$ws.Range("H122").Value = 26353014
$ws.Range("I122").Value = 14287568
$ws.Range("J122").Value = 39346572
$ws.Range("K122").Value = 42862704
$ws.Range("L122").Value = 118039716
$ws.Range("M122").Value = -42860254
$ws.Range("N122").Value = -118044616
$ws.Range("H137").Value = 11170351
$ws.Range("I137").Value = 928.9211
$ws.Range("J137").Value = 34750240
$ws.Range("K137").Value = 2786.7633
$ws.Range("L137").Value = 104250720
$ws.Range("M137").Value = -236.7633000000001
$ws.Range("N137").Value = -104255820
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 85188200
$ws.Range("I74").Value = 83334890
$ws.Range("J74").Value = 88894830
$ws.Range("K74").Value = 83334890
$ws.Range("L74").Value = 88894830
$ws.Range("M74").Value = -83334016
$ws.Range("N74").Value = -88896578
$ws.Range("H77").Value = 85188200
$ws.Range("I77").Value = 83334890
$ws.Range("J77").Value = 88894830
$ws.Range("K77").Value = 416674450
$ws.Range("L77").Value = 444474150
$ws.Range("M77").Value = -416670082
$ws.Range("N77").Value = -444482886
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 714958.6
$ws.Range("I107").Value = 833868.4399999999
$ws.Range("K107").Value = 833868.4399999999
$ws.Range("M107").Value = -831948.4399999999
$ws.Range("H108").Value = 34980
$ws.Range("J108").Value = 34980
$ws.Range("L108").Value = 34980
$ws.Range("N108").Value = -42660
$ws.Range("H134").Value = 12756211
$ws.Range("I134").Value = 13889774
$ws.Range("J134").Value = 5954833
$ws.Range("K134").Value = 41669322
$ws.Range("L134").Value = 17864499
$ws.Range("M134").Value = -41666787
$ws.Range("N134").Value = -17869569
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 964499.5
$ws.Range("I31").Value = 1316.3914
$ws.Range("J31").Value = 3296416.5
$ws.Range("K31").Value = 1316.3914
$ws.Range("L31").Value = 3296416.5
$ws.Range("M31").Value = -1021.3914
$ws.Range("N31").Value = -3297006.5
$ws.Range("H34").Value = 964499.5
$ws.Range("I34").Value = 1316.3914
$ws.Range("J34").Value = 3296416.5
$ws.Range("K34").Value = 1316.3914
$ws.Range("L34").Value = 3296416.5
$ws.Range("M34").Value = -1114.3914
$ws.Range("N34").Value = -3296820.5
$ws.Range("H58").Value = 970160.4
$ws.Range("I58").Value = 4844.56
$ws.Range("J58").Value = 2067110.1
$ws.Range("K58").Value = 4844.56
$ws.Range("L58").Value = 2067110.1
$ws.Range("M58").Value = -4641.56
$ws.Range("N58").Value = -2067516.1
$ws.Range("H60").Value = 11149
$ws.Range("J60").Value = 11149
$ws.Range("L60").Value = 11149
$ws.Range("N60").Value = -12171
$ws.Range("H132").Value = 1927.425
$ws.Range("I132").Value = 1527.0303
$ws.Range("K132").Value = 4581.090899999999
$ws.Range("M132").Value = -2051.090899999999
$ws.Range("H134").Value = 1057756.1
$ws.Range("I134").Value = 4943.129
$ws.Range("K134").Value = 14829.387
$ws.Range("M134").Value = -12294.387
$ws.Range("H136").Value = 970160.4
$ws.Range("I136").Value = 4844.56
$ws.Range("J136").Value = 2067110.1
$ws.Range("K136").Value = 14533.68
$ws.Range("L136").Value = 6201330.300000001
$ws.Range("M136").Value = -11983.68
$ws.Range("N136").Value = -6206430.300000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 88.14286
$ws.Range("J12").Value = 79.833336
$ws.Range("L12").Value = 239.500008
$ws.Range("N12").Value = -585.500008
$ws.Range("H123").Value = 500
$ws.Range("I123").Value = 500
$ws.Range("K123").Value = 1500
$ws.Range("M123").Value = 950
$ws.Range("H129").Value = 38097976
$ws.Range("I129").Value = 100001620
$ws.Range("J129").Value = 7146156
$ws.Range("K129").Value = 300004860
$ws.Range("L129").Value = 21438468
$ws.Range("M129").Value = -299999860
$ws.Range("N129").Value = -21448468
$ws.Range("H130").Value = 1418.2
$ws.Range("J130").Value = 1408.3
$ws.Range("L130").Value = 4224.9
$ws.Range("N130").Value = -14264.9
$ws.Range("H131").Value = 7057198
$ws.Range("I131").Value = 45545844
$ws.Range("J131").Value = 946.68335
$ws.Range("K131").Value = 136637532
$ws.Range("L131").Value = 2840.05005
$ws.Range("M131").Value = -136632492
$ws.Range("N131").Value = -12920.05005
$ws.Range("H132").Value = 3324.7407
$ws.Range("I132").Value = 2681.7778
$ws.Range("J132").Value = 3646.2222
$ws.Range("K132").Value = 24136.0002
$ws.Range("L132").Value = 32815.99980000001
$ws.Range("M132").Value = -21606.0002
$ws.Range("N132").Value = -37875.99980000001
$ws.Range("H133").Value = 3469.6553
$ws.Range("I133").Value = 4135
$ws.Range("J133").Value = 3000
$ws.Range("K133").Value = 12405
$ws.Range("L133").Value = 9000
$ws.Range("M133").Value = -7345
$ws.Range("N133").Value = -19120
$ws.Range("H134").Value = 2899
$ws.Range("I134").Value = 1623.75
$ws.Range("J134").Value = 8000
$ws.Range("K134").Value = 4871.25
$ws.Range("L134").Value = 24000
$ws.Range("M134").Value = 198.75
$ws.Range("N134").Value = -34140
$ws.Range("H136").Value = 2994.647
$ws.Range("I136").Value = 1847.1428
$ws.Range("J136").Value = 3797.9
$ws.Range("K136").Value = 5541.428400000001
$ws.Range("L136").Value = 11393.7
$ws.Range("M136").Value = -441.4284000000007
$ws.Range("N136").Value = -21593.7
$ws.Range("H137").Value = 2531.4167
$ws.Range("I137").Value = 2004
$ws.Range("J137").Value = 3410.4443
$ws.Range("K137").Value = 6012
$ws.Range("L137").Value = 10231.3329
$ws.Range("M137").Value = -912
$ws.Range("N137").Value = -20431.3329
$ws.Range("H138").Value = 111112870
$ws.Range("I138").Value = 142858350
$ws.Range("K138").Value = 428575050
$ws.Range("M138").Value = -428569910
$ws.Range("H139").Value = 56888.89
$ws.Range("I139").Value = 63625
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 190875
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = -185735
$ws.Range("N139").Value = -19280
$ws.Range("H140").Value = 3435.923
$ws.Range("I140").Value = 2453.4
$ws.Range("J140").Value = 4345.6665
$ws.Range("K140").Value = 7360.200000000001
$ws.Range("L140").Value = 13036.9995
$ws.Range("M140").Value = -2180.200000000001
$ws.Range("N140").Value = -23396.9995
$ws.Range("H141").Value = 2565
$ws.Range("I141").Value = 2343.6365
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 7030.9095
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -1850.9095
$ws.Range("N141").Value = -25360
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1857884.4
$ws.Range("I132").Value = 2600134.5
$ws.Range("J132").Value = 2259
$ws.Range("K132").Value = 7800403.5
$ws.Range("L132").Value = 6777
$ws.Range("M132").Value = -7797873.5
$ws.Range("N132").Value = -11837
$ws.Range("H136").Value = 1502523
$ws.Range("I136").Value = 1611201.2
$ws.Range("J136").Value = 2763
$ws.Range("K136").Value = 4833603.6
$ws.Range("L136").Value = 8289
$ws.Range("M136").Value = -4831053.6
$ws.Range("N136").Value = -13389
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 884499.25
$ws.Range("I132").Value = 3137.2424
$ws.Range("J132").Value = 2269496.8
$ws.Range("K132").Value = 9411.727200000001
$ws.Range("L132").Value = 6808490.399999999
$ws.Range("M132").Value = -6881.727200000001
$ws.Range("N132").Value = -6813550.399999999
$ws.Range("H136").Value = 3382.6765
$ws.Range("I136").Value = 2082.4866
$ws.Range("J136").Value = 4934.516
$ws.Range("K136").Value = 6247.459800000001
$ws.Range("L136").Value = 14803.548
$ws.Range("M136").Value = -3697.459800000001
$ws.Range("N136").Value = -19903.548
